# Trade #28 closed at 2026-02-17 13:19:15 - unknown UNKNOWN +0.000%
#
# Updates the Summary + Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly-closed trade #28 to both
# the "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.27   # Current Capital
$summary.Range("B4").Value = -0.74    # Total P&L $
$summary.Range("B5").Value = -0.53    # Total P&L %
$summary.Range("B6").Value = 28        # Total Trades
$summary.Range("B7").Value = 10        # Winning Trades
$summary.Range("B9").Value = 35.71     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.27      # Capital
$status.Range("D4").Value = 28         # Trades
$status.Range("E4").Value = -0.74     # P&L $
$status.Range("F4").Value = -0.73     # P&L %
$status.Range("G4").Value = 35.71      # Win Rate %

# ---------------------------------------------------------------------
# Append trade #28 to a trade-log sheet at row 29
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Range("A29").Value = 28
    $ws.Range("B29").NumberFormat = "@"
    $ws.Range("B29").Value = "2026-02-17"
    $ws.Range("B29").Style = "Normal"
    $ws.Range("C29").Value = "13:19:09"
    $ws.Range("D29").Value = "MarketMaking"
    $ws.Range("E29").Value = "DOWN"
    $ws.Range("F29").Value = 0.85
    $ws.Range("G29").Value = 0.89
    $ws.Range("H29").Value = "CLOSED"
    $ws.Range("I29").Value = 4.7059
    $ws.Range("J29").Value = 0.04
    $ws.Range("K29").Value = 99.27
    $ws.Range("L29").Value = 0
    $ws.Range("M29").Value = 0
    $ws.Range("N29").Value = 0.6
    $ws.Range("O29").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P29").Value = "early_exit"
    $ws.Range("Q29").Value = 0.13
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
